$d = $word.ActiveDocument

# Locate the paragraph containing "B2:" and then work on the (empty) paragraph
# that immediately follows it - that is the one that gets the new line of text.
$targetParagraph = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cur = $d.Paragraphs.Item($i)
    $curText = $cur.Range.Text.TrimEnd([char]13, [char]7)
    if ($curText -eq "B2:") {
        $targetParagraph = $d.Paragraphs.Item($i + 1)
        break
    }
}

$targetRange = $targetParagraph.Range

# Preserve the paragraph's existing identity attributes (w14:paraId, rsids, ...)
# if this runtime exposes them, but don't depend on it being there.
$existingXml = $targetRange.WordOpenXML
$paraAttrs = ""
if ($existingXml -match "<w:p( [^>]*?)/?>") {
    $paraAttrs = $matches[1]
}

# Rebuild the paragraph with a leading space run (normal formatting) followed
# by a 10pt ("Hello I am B2") run, and stamp the same 10pt size on the
# paragraph mark itself via pPr/rPr, matching what Word's Font Size box does
# when applied to the whole paragraph including its end-of-paragraph mark.
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"' + $paraAttrs + '><w:pPr><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Hello I am B2</w:t></w:r></w:p>'

$null = $targetRange.InsertXML($newParaXml)
